$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.693996120467564
$ws.Range("K2").Value = 0.691681993781983
$ws.Range("L2").Value = 0.754855071318372
$ws.Range("N2").Value = 0.622335501686223

$ws.Range("B3").Value = 0.641188431291775
$ws.Range("K3").Value = 0.550394111458431
$ws.Range("L3").Value = 0.756622372880961
$ws.Range("N3").Value = 0.57421952327369
